$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 34 previously described the 29-Maio exam date, duplicated into row 35.
# Now row 34 becomes "no class" (----) and row 35 becomes the actual exam date (29-Maio).
$ws.Range("C34").Value = "----"
$ws.Range("D34").Value = "Sem aula"
$ws.Range("E34").Value = "Sem aula"
$ws.Range("F34").Value = "Sem aula"
